$wb = $excel.ActiveWorkbook

function Remove-RowAndStaleHyperlinks($ws, $rowToDelete) {
    # Capture every existing hyperlink (address/display/target-range) before
    # we touch anything, because deleting the row shifts remaining rows up
    # and this engine's Hyperlinks.Delete() wipes the *whole* collection.
    $addrs = @()
    $disps = @()
    $refs = @()
    foreach ($hl in $ws.Hyperlinks) {
        $addrs += $hl.Address()
        $disps += $hl.TextToDisplay()
        $refs += $hl.Range.Address()
    }

    $ws.Rows.Item($rowToDelete).Delete()
    $ws.Hyperlinks.Delete()

    for ($i = 0; $i -lt $addrs.Count; $i++) {
        $ref = $refs[$i]
        # $ref looks like "$A$3" - pull the row number out
        $parts = $ref.Split('$')
        $colPart = $parts[1]
        $rowPart = [int]$parts[2]
        if ($rowPart -ne $rowToDelete) {
            $newRow = $rowPart
            if ($rowPart -gt $rowToDelete) {
                $newRow = $rowPart - 1
            }
            $cellref = $colPart + $newRow
            $ws.Hyperlinks.Add($ws.Range($cellref), $addrs[$i], "", "", $disps[$i])
        }
    }
}

# --- Sheet "Overview" ---
$wsOverview = $wb.Worksheets.Item("Overview")
Remove-RowAndStaleHyperlinks $wsOverview 3

# --- Sheet "zh-cn" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Remove-RowAndStaleHyperlinks $wsZhCn 3
$wsZhCn.Range("E2").Value = "2016-03-20 18:37:21"
$wsZhCn.Range("H2").Value = "2016-03-20 18:37:40"

# --- Sheet "de-de" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
Remove-RowAndStaleHyperlinks $wsDeDe 3
$wsDeDe.Range("E2").Value = "2016-03-20 18:37:24"
$wsDeDe.Range("H2").Value = "2016-03-20 18:37:45"
